$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 3: capacitor label gains C11
$ws.Range("C3").Value = "C1, C2, C3, C4, C11"

# 2. Row 4 (10k resistors): URL changed, hyperlink removed entirely (becomes plain text)
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$4') {
        $hl.Delete()
    }
}

# 3. Row 11 (30k resistors): URL changed, hyperlink removed entirely (becomes plain text)
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$11') {
        $hl.Delete()
    }
}

$ws.Range("B4").Value = "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/2230562/"
$ws.Range("B11").Value = "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6791263/"

# 4. Row 10 (1k resistors): URL changed, hyperlink target updated (still a live hyperlink)
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$10') {
        $hl.Address = "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/2230427/"
    }
}
$ws.Range("B10").Value = "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/2230427/"

# 5. New row 12: 10 pin female header (plain text, no hyperlink)
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A12").Value = "10 pin female header"
$ws.Range("B12").Value = "http://china.rs-online.com/web/p/pcb-sockets/7655745/"
$ws.Range("C12").Value = "J1, J2"
$ws.Range("D12").Value = 2

# 6. Selection moves to F4
[void]$ws.Range("F4").Select()
